# Fruta / hortaliza, semanal
# Insert a new weekly price-report row at row 162 (Femacal de La Calera,
# Arandano (blue)), pushing the existing rows 162-205 down to 163-206.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 162; this shifts rows
# 162..205 down to 163..206 (and grows the used range to A1:T206).
$ws.Rows("162").Insert()

# Populate the newly inserted row with this week's data point.
$ws.Range("A162").Value2 = 3
$ws.Range("B162").Value2 = "Femacal de La Calera"
$ws.Range("C162").Value2 = "Coquimbo"
$ws.Range("D162").Value2 = 44627
$ws.Range("E162").Value2 = 5
$ws.Range("F162").Value2 = "Fruta"
$ws.Range("G162").Value2 = 100101
$ws.Range("H162").Value2 = "Berries"
$ws.Range("I162").Value2 = 100101001
$ws.Range("J162").Value2 = "Arándano (blue)"
$ws.Range("K162").Value2 = "Sin especificar"
$ws.Range("L162").Value2 = "Primera"
$ws.Range("M162").Value2 = 50
$ws.Range("N162").Value2 = 4000
$ws.Range("O162").Value2 = 4000
$ws.Range("P162").Value2 = 4000
$ws.Range("Q162").Value2 = "$/bandeja 2 kilos"
$ws.Range("R162").Value2 = "Provincia de Linares"
$ws.Range("S162").Value2 = 2000
$ws.Range("T162").Value2 = 2
